# Update functions and Data Model (#50)
# Adds a new "Authorship Resource" column to Table1 / Sheet1, populated for
# both data rows, resizes the table + sheet dimension accordingly, and
# leaves the selection on L11 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table (and its AutoFilter) by one column: A1:J3 -> A1:K3.
$lo.Resize($ws.Range("A1:K3"))

# Header + data for the new "Authorship Resource" column.
$ws.Range("K1").Value = "Authorship Resource"
$ws.Range("K2").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("K3").Value = "Daniela Subotic, Noémi Villars-Amberg"

# Match the formatting (style) used by the rest of the header/data cells
# in the table (font/alignment style index), instead of the default style
# that plain value assignment would otherwise leave behind.
$ws.Range("J1:J3").Copy()
$ws.Range("K1:K3").PasteSpecial(-4122)

# Leave the selection where the author left it.
$ws.Range("L11").Select()
